$d = $word.ActiveDocument

# --- 1) Append two new paragraphs at the end of the document body ---
# First new paragraph: empty spacer paragraph (inherits the formatting
# of the preceding paragraph mark: not-bold, sz=32/32, en-IN).
$r1 = $d.Content
$r1.Collapse(0)
$r1.InsertParagraphAfter()

# Second new paragraph: contains the new requirement text.
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$r3 = $d.Content
$r3.Collapse(0)
$r3.InsertAfter("User after enter the Username and Password Login Button should be get Enabled")

# --- 2) Mark the built-in "Normal Table" style as a Quick Style ---
$tableStyle = $d.Styles("Normal Table")
$tableStyle.QuickStyle = $true
